# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, non-numeric-looking prices, volume %)
$ws.Range('D2').Value = '29.012.47'
$ws.Range('E2').Value = '  -4.20%  '
$ws.Range('D3').Value = '1.962.32'
$ws.Range('E3').Value = '  -6.15%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('E5').Value = '  -4.34%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  -6.01%  '
$ws.Range('E8').Value = '  -4.30%  '
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('E11').Value = '  -6.71%  '
$ws.Range('E12').Value = '  -7.55%  '
$ws.Range('D13').Value = '1.956.68'
$ws.Range('E13').Value = '  -7.78%  '
$ws.Range('E14').Value = '  -8.39%  '
$ws.Range('E15').Value = '  -6.67%  '
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('E17').Value = '  -5.32%  '
$ws.Range('E18').Value = '  -10.29%  '
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').Value = '  -9.29%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('E22').Value = '  -5.99%  '
$ws.Range('D23').Value = '29.046.46'
$ws.Range('E23').Value = '  -4.16%  '
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').Value = '2.241.95'
$ws.Range('E26').Value = '  -4.06%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E27').Value = '  -4.08%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E28').Value = '  -5.65%  '
$ws.Range('E29').Value = '  -10.36%  '
$ws.Range('E30').Value = '  -9.52%  '
$ws.Range('E31').Value = '  -5.28%  '
$ws.Range('E32').Value = '  -8.64%  '
$ws.Range('E33').Value = '  -6.46%  '
$ws.Range('E34').Value = '  -8.67%  '
$ws.Range('E35').Value = '  -8.31%  '
$ws.Range('E36').Value = '  -6.00%  '
$ws.Range('E37').Value = '  -7.75%  '
$ws.Range('E39').Value = '  -11.78%  '
$ws.Range('E40').Value = '  -6.65%  '
$ws.Range('E41').Value = '  -7.63%  '
$ws.Range('E42').Value = '  -9.44%  '
$ws.Range('E43').Value = '  -10.36%  '
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('E45').Value = '  -7.98%  '
$ws.Range('E47').Value = '  -8.87%  '
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('E50').Value = '  -4.50%  '
$ws.Range('E51').Value = '  -4.49%  '

# Price values that look numeric must stay literal text (e.g. "1.007"),
# so force Text format before assigning, then restore the default style
# so formatting matches the rest of the sheet.
$numericPriceCells = @{
    'D4' = '1.007'
    'D5' = '327.13'
    'D6' = '1.007'
    'D7' = '0.4991'
    'D8' = '0.4198'
    'D9' = '52.69'
    'D10' = '0.09150'
    'D11' = '1.096'
    'D12' = '22.83'
    'D14' = '7.841'
    'D15' = '6.425'
    'D16' = '1.008'
    'D17' = '0.00001098'
    'D18' = '91.19'
    'D19' = '0.06671'
    'D20' = '19.19'
    'D21' = '1.006'
    'D22' = '5.965'
    'D24' = '12.03'
    'D25' = '2.285'
    'D27' = '156.22'
    'D28' = '20.58'
    'D29' = '6.170'
    'D30' = '2.259'
    'D31' = '126.59'
    'D32' = '1.037'
    'D33' = '0.09835'
    'D34' = '1.524'
    'D35' = '5.760'
    'D36' = '3.677'
    'D37' = '0.02415'
    'D38' = '1.298'
    'D39' = '8.937'
    'D40' = '0.06313'
    'D41' = '0.6432'
    'D42' = '11.42'
    'D43' = '0.1981'
    'D44' = '1.006'
    'D45' = '0.6216'
    'D46' = '13.36'
    'D47' = '2.171'
    'D48' = '1.293'
    'D49' = '3.465'
    'D51' = '0.06936'
}
foreach ($cellRef in $numericPriceCells.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericPriceCells[$cellRef]
    $cell.Style = "Normal"
}
